$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48:107 down to 49:108
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new data point
$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 45036
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = 100112037
$ws.Cells.Item(48, 7).Value = "Cebollín"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 60
$ws.Cells.Item(48, 11).Value = 4000
$ws.Cells.Item(48, 12).Value = 4500
$ws.Cells.Item(48, 13).Value = 4250
$ws.Cells.Item(48, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 118
$ws.Cells.Item(48, 17).Value = 36
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the other date cells in column D
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(49, 4).NumberFormat
